$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.318.00"
$ws.Range("E2").Value = "  -2.68%  "

$ws.Range("D3").Value = "1.569.79"
$ws.Range("E3").Value = "  -3.69%  "

$ws.Range("E4").Value = "  -0.43%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "207.85"
$ws.Range("E5").Value = "  -2.97%  "

$ws.Range("E6").Value = "  -0.42%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.476"
$ws.Range("E7").Value = "  -5.19%  "

$ws.Range("E8").Value = "  -2.10%  "

$ws.Range("E9").Value = "  -1.82%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "17.92"
$ws.Range("E10").Value = "  -1.88%  "

$ws.Range("E11").Value = "  -0.90%  "

$ws.Range("D12").Value = "1.789.11"
$ws.Range("E12").Value = "  -3.62%  "

$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.574.30"
$ws.Range("E13").Value = "  -5.69%  "

$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.03"
$ws.Range("E14").Value = "  -3.21%  "

$ws.Range("E15").Value = "  -3.07%  "

$ws.Range("D16").Value = "25.315.46"
$ws.Range("E16").Value = "  -2.58%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "59.69"
$ws.Range("E17").Value = "  -2.56%  "

$ws.Range("E18").Value = "  -4.14%  "

$ws.Range("E19").Value = "  -0.40%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "185.65"
$ws.Range("E20").Value = "  -2.18%  "

$ws.Range("E21").Value = "  -2.06%  "

$ws.Range("E22").Value = "  -2.45%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.88"
$ws.Range("E23").Value = "  -2.86%  "

$ws.Range("E24").Value = "  -1.96%  "

$ws.Range("E25").Value = "  -0.43%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "141.08"
$ws.Range("E26").Value = "  -1.80%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.66"
$ws.Range("E27").Value = "  -6.98%  "

$ws.Range("B28").Value = "EthereumClassic"
$ws.Range("C28").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "14.88"
$ws.Range("E28").Value = "  -1.47%  "

$ws.Range("B29").Value = "Cosmos"
$ws.Range("C29").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.46"
$ws.Range("E29").Value = "  -3.87%  "

$ws.Range("E30").Value = "  -6.02%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0462"
$ws.Range("E31").Value = "  -3.73%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.05"
$ws.Range("E32").Value = "  -2.35%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.00"
$ws.Range("E33").Value = "  -3.88%  "

$ws.Range("E34").Value = "  -1.19%  "

$ws.Range("E35").Value = "  -4.46%  "

$ws.Range("D36").Value = "1.090.43"
$ws.Range("E36").Value = "  -3.56%  "

$ws.Range("E37").Value = "  -0.72%  "

$ws.Range("B38").Value = "MXToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.32"
$ws.Range("E38").Value = "  -4.89%  "

$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0150"
$ws.Range("E39").Value = "  -2.37%  "

$ws.Range("E40").Value = "  -3.66%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.774"
$ws.Range("E41").Value = "  -8.76%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.767"
$ws.Range("E42").Value = "  -0.79%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "92.47"
$ws.Range("E43").Value = "  -5.63%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.06"
$ws.Range("E44").Value = "  -2.70%  "

$ws.Range("D45").Value = "1.703.01"
$ws.Range("E45").Value = "  -3.63%  "

$ws.Range("D46").Value = "0.0₆0111"
$ws.Range("E46").Value = "  -2.71%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "52.86"
$ws.Range("E47").Value = "  -3.17%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0505"
$ws.Range("E48").Value = "  -3.96%  "

$ws.Range("E49").Value = "  -3.94%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.407"
$ws.Range("E50").Value = "  -1.77%  "

$ws.Range("E51").Value = "  -0.47%  "
